# Ajuste nuevas ligas, porteros y metricas
# Rewrites the Competencia / Ponderacion_Competencia table (A1:B41)
# with the updated league list and re-ordered weights.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$competencias = @(
    "England Premier League",
    "Europe Champions League",
    "Germany Bundesliga",
    "Italy Serie A",
    "Spain La Liga",
    "France Ligue 1",
    "Europa League",
    "Brasileirao",
    "Liga Portugal",
    "Liga Profesional Argentina",
    "Eredivisie",
    "MLS",
    "Jupiler Pro League",
    "Spanish Segunda Division",
    "Austrian Bundesliga",
    "Brack Super League",
    "Croatia Prva HNL",
    "Danish Superligaen",
    "England Championship",
    "Italian Serie B",
    "Turkey Super Lig",
    "UEFA Under 17 Championship",
    "UEFA Under 19 Championship",
    "UEFA Under 21 Championship",
    "Colombia Superliga",
    "Scotland Premiership",
    "CONMEBOL Libertadores U20",
    "CONMEBOL U17",
    "CONMEBOL U20",
    "Ecuador Liga Pro",
    "Chile Primera",
    "Belgian Challenger Pro League",
    "Bundesliga 2",
    "Dutch Eerste Divisie",
    "French Ligue 2",
    "Norwegian Eliteserien",
    "Polish Ekstraklasa",
    "Russia Premier League",
    "Serbian Super Liga",
    "Swedish Allsvenskan"
)
$ponderaciones = @(
    1.9,
    1.9,
    1.9,
    1.9,
    1.9,
    1.7,
    1.5,
    1.4,
    1.3,
    1.3,
    1.2,
    1.2,
    1.1,
    1.1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    0.9,
    0.9,
    1,
    1,
    1,
    0.9,
    0.8,
    0.8,
    0.8,
    0.8,
    0.8,
    0.8,
    0.7,
    0.7,
    0.7,
    0.7
)

for ($i = 0; $i -lt $competencias.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $competencias[$i]
    $ws.Cells.Item($row, 2).Value = $ponderaciones[$i]
}

# Restore the view state recorded in the edited workbook:
# scrolled so row 24 is at the top, with A41 selected.
$excel.ActiveWindow.ScrollRow = 24
$ws.Range("A41").Select()
